$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The tracked data lives in rows 2..184 (header row 1), filtered on
# D (Difficulty) = "Medium" and E (Finished) = "N".
$firstDataRow = 2
$lastDataRow = 184

# Snapshot every data row's current Hidden state so that reapplying the
# AutoFilter below (needed to refresh its saved metadata) cannot change the
# visibility of any row other than the one we are intentionally completing.
$hiddenState = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $hiddenState[$r] = $ws.Rows($r).Hidden
}

# 95. Unique Binary Search Trees II -> mark it Finished.
$ws.Range("E96").Value = "Y"
$hiddenState[96] = $true

# Reapply the AutoFilter on the "Finished" column (field 2 of the D:E
# filter range) with its existing criterion so the persisted filter
# definition is refreshed.
$ws.AutoFilter.Range.AutoFilter(2, @("N"), 7)

# Restore every row's visibility to what it was beforehand, except for the
# row we just marked Finished, which should now be hidden by the filter.
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ws.Rows($r).Hidden = $hiddenState[$r]
}

# Move the active selection/cursor to G104.
$ws.Range("G104").Select()
